# Refresh crypto price/volume figures in the worksheet to match the
# latest scrape (GitHub Actions data-refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.217.36'
$ws.Range("E2").Value = '  -4.40%  '
$ws.Range("D3").Value = '2.933.14'
$ws.Range("E3").Value = '  -7.16%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''479.32'
$ws.Range("E5").Value = '  -8.84%  '
$ws.Range("D6").Value = '''129.18'
$ws.Range("E6").Value = '  -3.18%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '2.933.54'
$ws.Range("E8").Value = '  -7.12%  '
$ws.Range("E9").Value = '  -8.65%  '
$ws.Range("D10").Value = '''6.89'
$ws.Range("E10").Value = '  -5.58%  '
$ws.Range("D11").Value = '''0.0994'
$ws.Range("E11").Value = '  -10.40%  '
$ws.Range("E12").Value = '  -11.99%  '
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").Value = '3.432.93'
$ws.Range("E14").Value = '  -7.36%  '
$ws.Range("D15").Value = '''24.02'
$ws.Range("E15").Value = '  -7.09%  '
$ws.Range("D16").Value = '55.175.43'
$ws.Range("E16").Value = '  -4.48%  '
$ws.Range("D17").Value = '2.934.42'
$ws.Range("E17").Value = '  -7.35%  '
$ws.Range("D18").Value = '''0.0000138'
$ws.Range("E18").Value = '  -9.82%  '
$ws.Range("D19").Value = '''5.52'
$ws.Range("E19").Value = '  -5.19%  '
$ws.Range("D20").Value = '''11.80'
$ws.Range("E20").Value = '  -9.72%  '
$ws.Range("D21").Value = '''7.34'
$ws.Range("E21").Value = '  -9.04%  '
$ws.Range("D22").Value = '''307.77'
$ws.Range("E22").Value = '  -11.16%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  -11.32%  '
$ws.Range("D25").Value = '''59.46'
$ws.Range("E25").Value = '  -14.59%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '''0.156'
$ws.Range("E27").Value = '  -6.31%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '0.0₃0833'
$ws.Range("E29").Value = '  -13.19%  '
$ws.Range("D30").Value = '''6.44'
$ws.Range("E30").Value = '  -6.23%  '
$ws.Range("D31").Value = '''1.16'
$ws.Range("E31").Value = '  -4.84%  '
$ws.Range("E32").Value = '  -7.63%  '
$ws.Range("E33").Value = '  -12.09%  '
$ws.Range("D34").Value = '''19.07'
$ws.Range("E34").Value = '  -12.33%  '
$ws.Range("D35").Value = '''146.52'
$ws.Range("E35").Value = '  -8.47%  '
$ws.Range("D36").Value = '''4.29'
$ws.Range("E36").Value = '  -12.53%  '
$ws.Range("D37").Value = '''5.56'
$ws.Range("E37").Value = '  -11.12%  '
$ws.Range("E38").Value = '  -11.07%  '
$ws.Range("D39").Value = '''23.36'
$ws.Range("E39").Value = '  -9.90%  '
$ws.Range("D40").Value = '''0.0637'
$ws.Range("E40").Value = '  -8.58%  '
$ws.Range("D41").Value = '2.960.55'
$ws.Range("E41").Value = '  -7.18%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").Value = '''35.78'
$ws.Range("E43").Value = '  -11.86%  '
$ws.Range("D44").Value = '''0.984'
$ws.Range("E44").Value = '  -9.17%  '
$ws.Range("E45").Value = '  -11.00%  '
$ws.Range("D46").Value = '''1.35'
$ws.Range("E46").Value = '  -7.82%  '
$ws.Range("E47").Value = '  -11.92%  '
$ws.Range("D48").Value = '2.100.44'
$ws.Range("E48").Value = '  -7.50%  '
$ws.Range("D49").Value = '''0.0225'
$ws.Range("E49").Value = '  -4.79%  '
$ws.Range("D50").Value = '''18.59'
$ws.Range("E50").Value = '  -9.38%  '
$ws.Range("E51").Value = '  -11.33%  '
